$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove all existing hyperlinks (and their relationships) so we can
#    rebuild the sheet content/relationships cleanly and in the right order.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 2. Preserve formatting for the row/cells that are new (row 8 did not exist
#    before; column F did not exist before) by copying the existing look of
#    analogous cells before we touch values/hyperlinks.
# ---------------------------------------------------------------------------
# A8 should look like the other index cells in column A (bold, bordered).
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(8, 1).PasteSpecial(-4122) | Out-Null

# B/D columns (2-8) should carry the "Hyperlink" cell style, same as before.
for ($i = 2; $i -le 8; $i++) {
    $ws.Cells.Item($i, 2).Style = "Hyperlink"
    $ws.Cells.Item($i, 4).Style = "Hyperlink"
}

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Add the new "MOP_DEF" header column (F1).
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 6).Value = "MOP_DEF"

# ---------------------------------------------------------------------------
# 4. Rewrite the data rows (2-8). Row 2 is a brand-new row for CHEBI_23367
#    (with its textual definition in column F); the previously-existing
#    rows are reflowed below it (note CHEBI_52214 / CHEBI_39141 swap order),
#    and the final BAO_0003043 / CHEBI_23367 row moves from row 7 to row 8
#    and also gains the CHEBI_23367 definition in column F.
# ---------------------------------------------------------------------------

$definition = "['Any constitutionally or isotopically distinct atom, molecule, ion, ion pair, radical, radical ion, complex, conformer etc., identifiable as a separately distinguishable entity. [IUPAC]']"
$noDefinition = "[]"

$data = @(
    @{ A = 0; B = "http://purl.obolibrary.org/obo/CHEBI_23367"; C = "{'iri': 'http://purl.obolibrary.org/obo/CHEBI_23367'}"; D = "http://purl.obolibrary.org/obo/CHEBI_23367"; E = "{'iri': 'http://purl.obolibrary.org/obo/CHEBI_23367'}"; F = $definition },
    @{ A = 1; B = "http://purl.obolibrary.org/obo/CHEBI_39141"; C = "{'iri': 'http://purl.obolibrary.org/obo/CHEBI_39141'}"; D = "http://purl.obolibrary.org/obo/CHEBI_39141"; E = "{'iri': 'http://purl.obolibrary.org/obo/CHEBI_39141'}"; F = $noDefinition },
    @{ A = 2; B = "http://purl.obolibrary.org/obo/CHEBI_52214"; C = "{'iri': 'http://purl.obolibrary.org/obo/CHEBI_52214'}"; D = "http://purl.obolibrary.org/obo/CHEBI_52214"; E = "{'iri': 'http://purl.obolibrary.org/obo/CHEBI_52214'}"; F = $noDefinition },
    @{ A = 3; B = "http://purl.obolibrary.org/obo/CHEBI_39142"; C = "{'iri': 'http://purl.obolibrary.org/obo/CHEBI_39142'}"; D = "http://purl.obolibrary.org/obo/CHEBI_39142"; E = "{'iri': 'http://purl.obolibrary.org/obo/CHEBI_39142'}"; F = $noDefinition },
    @{ A = 4; B = "http://purl.obolibrary.org/obo/CHEBI_39143"; C = "{'iri': 'http://purl.obolibrary.org/obo/CHEBI_39143'}"; D = "http://purl.obolibrary.org/obo/CHEBI_39143"; E = "{'iri': 'http://purl.obolibrary.org/obo/CHEBI_39143'}"; F = $noDefinition },
    @{ A = 5; B = "http://purl.obolibrary.org/obo/CHEBI_39144"; C = "{'iri': 'http://purl.obolibrary.org/obo/CHEBI_39144'}"; D = "http://purl.obolibrary.org/obo/CHEBI_39144"; E = "{'iri': 'http://purl.obolibrary.org/obo/CHEBI_39144'}"; F = $noDefinition },
    @{ A = 6; B = "http://www.bioassayontology.org/bao#BAO_0003043"; C = "{'label': 'molecular entity', 'prefLabel': None, 'altLabel': None, 'name': 'BAO_0003043'}"; D = "http://purl.obolibrary.org/obo/CHEBI_23367"; E = "{'label': 'molecular entity', 'prefLabel': 'molecular entity'}"; F = $definition }
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData.A
    $ws.Cells.Item($row, 2).Value = $rowData.B
    $ws.Cells.Item($row, 3).Value = $rowData.C
    $ws.Cells.Item($row, 4).Value = $rowData.D
    $ws.Cells.Item($row, 5).Value = $rowData.E
    $ws.Cells.Item($row, 6).Value = $rowData.F
}

# ---------------------------------------------------------------------------
# 5. Re-create the hyperlinks for columns B and D on rows 2-8. Row 8's B
#    hyperlink carries the in-document "BAO_0003043" bookmark/location like
#    the original row 7 hyperlink did.
# ---------------------------------------------------------------------------

$links = @(
    @{ Cell = "B2"; Url = "http://purl.obolibrary.org/obo/CHEBI_23367"; Loc = "" },
    @{ Cell = "D2"; Url = "http://purl.obolibrary.org/obo/CHEBI_23367"; Loc = "" },
    @{ Cell = "B3"; Url = "http://purl.obolibrary.org/obo/CHEBI_39141"; Loc = "" },
    @{ Cell = "D3"; Url = "http://purl.obolibrary.org/obo/CHEBI_39141"; Loc = "" },
    @{ Cell = "B4"; Url = "http://purl.obolibrary.org/obo/CHEBI_52214"; Loc = "" },
    @{ Cell = "D4"; Url = "http://purl.obolibrary.org/obo/CHEBI_52214"; Loc = "" },
    @{ Cell = "B5"; Url = "http://purl.obolibrary.org/obo/CHEBI_39142"; Loc = "" },
    @{ Cell = "D5"; Url = "http://purl.obolibrary.org/obo/CHEBI_39142"; Loc = "" },
    @{ Cell = "B6"; Url = "http://purl.obolibrary.org/obo/CHEBI_39143"; Loc = "" },
    @{ Cell = "D6"; Url = "http://purl.obolibrary.org/obo/CHEBI_39143"; Loc = "" },
    @{ Cell = "B7"; Url = "http://purl.obolibrary.org/obo/CHEBI_39144"; Loc = "" },
    @{ Cell = "D7"; Url = "http://purl.obolibrary.org/obo/CHEBI_39144"; Loc = "" },
    @{ Cell = "B8"; Url = "http://www.bioassayontology.org/bao"; Loc = "BAO_0003043" },
    @{ Cell = "D8"; Url = "http://purl.obolibrary.org/obo/CHEBI_23367"; Loc = "" }
)

foreach ($link in $links) {
    $ws.Hyperlinks.Add($ws.Range($link.Cell), $link.Url, $link.Loc)
}

# ---------------------------------------------------------------------------
# 6. Adding hyperlinks can perturb the cell style reference (Excel tends to
#    materialize a fresh style record); reapply the "Hyperlink" named style
#    so B/D columns keep referencing the original shared style, matching
#    the look of the untouched rows.
# ---------------------------------------------------------------------------
for ($i = 2; $i -le 8; $i++) {
    $ws.Cells.Item($i, 2).Style = "Hyperlink"
    $ws.Cells.Item($i, 4).Style = "Hyperlink"
}
